# Enter the NFL tip scores (predictions) for rows 57-60 (LAR vs NO, MIN vs
# DET, NE at DEN, KC vs LV) in columns E:H (Jasmin/Franz/David/Thomas).
# The assignment order below reproduces the exact shared-string insertion
# order recorded for this edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tipps")

$ws.Range("E58").Value = "DET 24-21"
$ws.Range("H59").Value = "DEN 20-14"
$ws.Range("H60").Value = "KC  27-16"
$ws.Range("E59").Value = "DEN 24-14"
$ws.Range("G57").Value = "LAR 24-23"
$ws.Range("G58").Value = "DET 27-24"
$ws.Range("G59").Value = "DEN 21-13"
$ws.Range("F58").Value = "DET 27-20"
$ws.Range("F59").Value = "DEN 20-13"
$ws.Range("F60").Value = "KC  17-10"
$ws.Range("E57").Value = "LAR 27-24"
$ws.Range("F57").Value = "LAR 23-20"
$ws.Range("H57").Value = "LAR 24-21"
$ws.Range("H58").Value = "DET 24-21"
$ws.Range("E60").Value = "KC  27-21"
$ws.Range("G60").Value = "KC  27-17"

$ws.Range("D57").Select()

$wb.Save()
